$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The small "Year / Month / National unemployment rate" table lives in
# columns X:Z, starting at row 6 (row 4 = headers, row 5 = spacer).
# A new, more recent monthly reading (2025 / Mar. / 2.2221963213599998) is
# being added at the top of the table, so every existing row shifts down
# by one (old row 19 -> row 20, ..., old row 6 -> row 7), and the new data
# is written into row 6.

$firstRow = 6
$lastRow  = 19

# Shift existing rows down by one, starting from the bottom so we never
# overwrite a source row before it has been copied.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $srcRange = "X" + $r + ":Z" + $r
    $dstRange = "X" + ($r + 1) + ":Z" + ($r + 1)
    $src = $ws.Range($srcRange)
    $dst = $ws.Range($dstRange)

    # Make sure the destination cells already exist so that pasting formats
    # lands on the exact same style as the source cells.
    $dst.Value = 0

    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats

    $src.Copy()
    $dst.PasteSpecial(-4163)  # xlPasteValues
}

# Write the new first row. Its formatting matches the style two rows below
# (the same alternating banding the table already used), so copy that.
$newDataRange  = $ws.Range("X6:Z6")
$bandSourceRange = $ws.Range("X8:Z8")

$newDataRange.Value = 0
$bandSourceRange.Copy()
$newDataRange.PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("X6").Value = 2025
$ws.Range("Y6").Value = "Mar."
$ws.Range("Z6").Value = 2.2221963213599998

$excel.CutCopyMode = 0
